# "Add filter and add_col scripts"
# Adds a "qty" column (D) to the expenses sheet and extends the sample
# data from 10 rows to 25 rows (rows 2-26), refreshing column C (price)
# values along the way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column D
$ws.Cells.Item(1, 4).Value = "qty"

# Full data block for rows 2-26: person_number, fruit, price, qty
$data = @(
    @(1,   "apple",      65, 8),
    @(1,   "strawberry", 71, 3),
    @(3,   "peer",       14, 1),
    @(4,   "apple",      71, 4),
    @(5,   "strawberry", 33, 3),
    @(6,   "lemon",      97, 6),
    @(7,   "kiwi",        1, 7),
    @(7,   "apple",      26, 3),
    @(1,   "peer",        5, 3),
    @(102, "orange",     47, 9),
    @(3,   "apple",      29, 10),
    @(4,   "strawberry",  6, 2),
    @(5,   "peer",       76, 9),
    @(6,   "apple",      18, 1),
    @(7,   "strawberry", 90, 4),
    @(3,   "lemon",      25, 2),
    @(4,   "kiwi",       62, 9),
    @(5,   "apple",      95, 9),
    @(6,   "peer",       96, 5),
    @(7,   "orange",     99, 9),
    @(7,   "apple",      69, 1),
    @(1,   "peer",       81, 3),
    @(102, "orange",     12, 5),
    @(3,   "apple",      18, 3),
    @(103, "strawberry", 35, 5)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Match the saved selection state from the edit
$ws.Range("C25").Select()
